$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tag_dict")

# Column C (data_type): cells that said "String" become "object"
$ws.Range("C4").Value = "object"
$ws.Range("C5").Value = "object"
$ws.Range("C6").Value = "object"
$ws.Range("C7").Value = "object"
$ws.Range("C8").Value = "object"
$ws.Range("C9").Value = "object"
$ws.Range("C10").Value = "object"
$ws.Range("C11").Value = "object"
$ws.Range("C12").Value = "object"
$ws.Range("C16").Value = "object"
$ws.Range("C23").Value = "object"
$ws.Range("C27").Value = "object"
$ws.Range("C28").Value = "object"
$ws.Range("C29").Value = "object"

# Column D (data_type_new): normalized to lowercase pandas dtype names
$ws.Range("D2").Value = "object"
$ws.Range("D3").Value = "int64"
$ws.Range("D4").Value = "object"
$ws.Range("D5").Value = "object"
$ws.Range("D6").Value = "object"
$ws.Range("D7").Value = "object"
$ws.Range("D8").Value = "object"
$ws.Range("D9").Value = "object"
$ws.Range("D10").Value = "object"
$ws.Range("D11").Value = "object"
$ws.Range("D12").Value = "object"
$ws.Range("D13").Value = "object"
$ws.Range("D14").Value = "int64"
$ws.Range("D15").Value = "int64"
$ws.Range("D16").Value = "object"
$ws.Range("D17").Value = "float64"
$ws.Range("D18").Value = "int64"
$ws.Range("D19").Value = "float64"
$ws.Range("D20").Value = "int64"
$ws.Range("D21").Value = "int64"
$ws.Range("D22").Value = "float64"
$ws.Range("D23").Value = "object"
$ws.Range("D24").Value = "int64"
$ws.Range("D25").Value = "float64"
$ws.Range("D26").Value = "int64"
$ws.Range("D27").Value = "object"
$ws.Range("D28").Value = "object"
$ws.Range("D29").Value = "object"
$ws.Range("D30").Value = "float64"
$ws.Range("D31").Value = "float64"
$ws.Range("D32").Value = "int64"
$ws.Range("D33").Value = "float64"
$ws.Range("D34").Value = "float64"
$ws.Range("D35").Value = "int64"
$ws.Range("D36").Value = "float64"
$ws.Range("D37").Value = "float64"
$ws.Range("D38").Value = "int64"
$ws.Range("D39").Value = "float64"
$ws.Range("D40").Value = "int64"
